$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to be inserted right after the header row (row 1),
# pushing the existing data down.
$newRows = @(
    @(0.0018325957935303, 0.1020144969224929, 0.0574213340878486),
    @(-0.1429424732923507, 0.204487144947052, 0.0444404482841491),
    @(-0.2125810980796814, 0.4230241775512695, 0.0167987942695617),
    @(-0.0858265683054924, 0.5377141237258911, -0.1299615800380706),
    @(-0.1346957832574844, 0.5236642360687256, 0.1327104717493057)
)

# Insert 5 new blank rows starting at row 2, shifting existing data down by 5.
$insertRange = $ws.Range("A2:C6")
$insertRange.EntireRow.Insert()

# Populate the newly inserted rows with the new data values.
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = 2 + $i
    $ws.Cells.Item($rowNum, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($rowNum, 3).Value = $newRows[$i][2]
}

# The oldest 3 rows of data (originally rows 17-19, now shifted to rows 22-24)
# are dropped so the dataset keeps a fixed 20-row window (rows 2-21).
$ws.Range("A22:C24").EntireRow.Delete()
